$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.190.87'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '1.659.41'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5165'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.26%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2641'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07760'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.477'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '1.649.81'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '1.886.90'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5458'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").Value = '0.0₅8125'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").Value = '26.218.00'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.613'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.990'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.96%  '
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1222'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.278'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.438'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05941'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.73%  '
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.544'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.268'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.585'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9619'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.52%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.769'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5668'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.045'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8540'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").Value = '1.011.09'
$ws.Range("E43").Value = '  -7.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '1.800.42'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.034'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05164'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4200'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '
